# This script reproduces the refreshed "Relatorio" currency-quote log:
#  - corrects the Hora/Data/Valor(R$) values for the existing 80 data rows
#    (the scraper re-ran and the later samples shifted into the earlier rows)
#  - appends 4 new rows (82-85) for the newest sample batch
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds comma-decimal text such as "4,91". A bare assignment of such a
# string gets misread as a *number* (the comma is treated as a thousands
# separator), so every such cell is switched to Text format first (cell by cell -
# a multi-area Range.NumberFormat assignment here only affects the first area),
# the value is assigned, and afterwards the cell is put back on the workbooks
# default "Normal" style so no stray formatting is left behind.
$cCells = @(6, 7, 9, 11, 19, 27, 31, 38, 39, 40, 41, 43, 54, 55, 56, 57, 59, 66, 67, 68, 69, 71, 82, 83, 84, 85)
foreach ($r in $cCells) { $ws.Cells.Item($r, 3).NumberFormat = "@" }

# --- Column C ("Valor(R$)") corrections for existing rows ---
$ws.Cells.Item(6,3).Value = '4,91'
$ws.Cells.Item(7,3).Value = '5,45'
$ws.Cells.Item(9,3).Value = '0,71'
$ws.Cells.Item(11,3).Value = '5,44'
$ws.Cells.Item(19,3).Value = '5,45'
$ws.Cells.Item(27,3).Value = '5,44'
$ws.Cells.Item(31,3).Value = '5,45'
$ws.Cells.Item(38,3).Value = '4,93'
$ws.Cells.Item(39,3).Value = '5,44'
$ws.Cells.Item(40,3).Value = '0,037'
$ws.Cells.Item(41,3).Value = '0,72'
$ws.Cells.Item(43,3).Value = '5,48'
$ws.Cells.Item(54,3).Value = '5,03'
$ws.Cells.Item(55,3).Value = '5,49'
$ws.Cells.Item(56,3).Value = '0,038'
$ws.Cells.Item(57,3).Value = '0,73'
$ws.Cells.Item(59,3).Value = '5,45'
$ws.Cells.Item(66,3).Value = '4,93'
$ws.Cells.Item(67,3).Value = '5,44'
$ws.Cells.Item(68,3).Value = '0,037'
$ws.Cells.Item(69,3).Value = '0,72'
$ws.Cells.Item(71,3).Value = '5,48'

# --- Column D ("Hora") and column E ("Data") corrections for existing rows ---
$ws.Cells.Item(2,4).Value = '12:31'
$ws.Cells.Item(3,4).Value = '12:31'
$ws.Cells.Item(4,4).Value = '12:31'
$ws.Cells.Item(5,4).Value = '12:31'
$ws.Cells.Item(6,4).Value = '12:29'
$ws.Cells.Item(6,5).Value = ' sábado, 15 de abril de 2023 '
$ws.Cells.Item(7,4).Value = '12:29'
$ws.Cells.Item(7,5).Value = ' sábado, 15 de abril de 2023 '
$ws.Cells.Item(8,4).Value = '12:29'
$ws.Cells.Item(8,5).Value = ' sábado, 15 de abril de 2023 '
$ws.Cells.Item(9,4).Value = '12:29'
$ws.Cells.Item(9,5).Value = ' sábado, 15 de abril de 2023 '
$ws.Cells.Item(10,4).Value = '18:31'
$ws.Cells.Item(11,4).Value = '18:31'
$ws.Cells.Item(12,4).Value = '18:31'
$ws.Cells.Item(13,4).Value = '18:31'
$ws.Cells.Item(14,4).Value = '17:50'
$ws.Cells.Item(15,4).Value = '17:50'
$ws.Cells.Item(16,4).Value = '17:50'
$ws.Cells.Item(17,4).Value = '17:50'
$ws.Cells.Item(18,4).Value = '17:48'
$ws.Cells.Item(19,4).Value = '17:48'
$ws.Cells.Item(20,4).Value = '17:48'
$ws.Cells.Item(21,4).Value = '17:48'
$ws.Cells.Item(22,4).Value = '17:47'
$ws.Cells.Item(23,4).Value = '17:47'
$ws.Cells.Item(24,4).Value = '17:47'
$ws.Cells.Item(25,4).Value = '17:47'
$ws.Cells.Item(26,4).Value = '17:45'
$ws.Cells.Item(27,4).Value = '17:45'
$ws.Cells.Item(28,4).Value = '17:45'
$ws.Cells.Item(29,4).Value = '17:45'
$ws.Cells.Item(30,4).Value = '17:38'
$ws.Cells.Item(31,4).Value = '17:38'
$ws.Cells.Item(32,4).Value = '17:38'
$ws.Cells.Item(33,4).Value = '17:38'
$ws.Cells.Item(34,4).Value = '17:32'
$ws.Cells.Item(35,4).Value = '17:32'
$ws.Cells.Item(36,4).Value = '17:32'
$ws.Cells.Item(37,4).Value = '17:32'
$ws.Cells.Item(38,4).Value = '17:29'
$ws.Cells.Item(38,5).Value = ' quinta-feira, 13 de abril de 2023 '
$ws.Cells.Item(39,4).Value = '17:29'
$ws.Cells.Item(39,5).Value = ' quinta-feira, 13 de abril de 2023 '
$ws.Cells.Item(40,4).Value = '17:29'
$ws.Cells.Item(40,5).Value = ' quinta-feira, 13 de abril de 2023 '
$ws.Cells.Item(41,4).Value = '17:29'
$ws.Cells.Item(41,5).Value = ' quinta-feira, 13 de abril de 2023 '
$ws.Cells.Item(42,4).Value = '22:32'
$ws.Cells.Item(43,4).Value = '22:32'
$ws.Cells.Item(44,4).Value = '22:32'
$ws.Cells.Item(45,4).Value = '22:32'
$ws.Cells.Item(46,4).Value = '21:21'
$ws.Cells.Item(47,4).Value = '21:21'
$ws.Cells.Item(48,4).Value = '21:21'
$ws.Cells.Item(49,4).Value = '21:21'
$ws.Cells.Item(54,4).Value = '21:20'
$ws.Cells.Item(54,5).Value = ' quarta-feira, 5 de abril de 2023 '
$ws.Cells.Item(55,4).Value = '21:20'
$ws.Cells.Item(55,5).Value = ' quarta-feira, 5 de abril de 2023 '
$ws.Cells.Item(56,4).Value = '21:20'
$ws.Cells.Item(56,5).Value = ' quarta-feira, 5 de abril de 2023 '
$ws.Cells.Item(57,4).Value = '21:20'
$ws.Cells.Item(57,5).Value = ' quarta-feira, 5 de abril de 2023 '
$ws.Cells.Item(58,4).Value = '17:38'
$ws.Cells.Item(59,4).Value = '17:38'
$ws.Cells.Item(60,4).Value = '17:38'
$ws.Cells.Item(61,4).Value = '17:38'
$ws.Cells.Item(62,4).Value = '17:32'
$ws.Cells.Item(63,4).Value = '17:32'
$ws.Cells.Item(64,4).Value = '17:32'
$ws.Cells.Item(65,4).Value = '17:32'
$ws.Cells.Item(66,4).Value = '17:29'
$ws.Cells.Item(66,5).Value = ' quinta-feira, 13 de abril de 2023 '
$ws.Cells.Item(67,4).Value = '17:29'
$ws.Cells.Item(67,5).Value = ' quinta-feira, 13 de abril de 2023 '
$ws.Cells.Item(68,4).Value = '17:29'
$ws.Cells.Item(68,5).Value = ' quinta-feira, 13 de abril de 2023 '
$ws.Cells.Item(69,4).Value = '17:29'
$ws.Cells.Item(69,5).Value = ' quinta-feira, 13 de abril de 2023 '
$ws.Cells.Item(70,4).Value = '22:32'
$ws.Cells.Item(71,4).Value = '22:32'
$ws.Cells.Item(72,4).Value = '22:32'
$ws.Cells.Item(73,4).Value = '22:32'
$ws.Cells.Item(74,4).Value = '21:21'
$ws.Cells.Item(75,4).Value = '21:21'
$ws.Cells.Item(76,4).Value = '21:21'
$ws.Cells.Item(77,4).Value = '21:21'

# --- Append the 4 new data rows (82-85), copying the style block (border +
#     bold + centered "A" column) from the row directly above ---
$ws.Cells.Item(81,1).Copy()
$ws.Cells.Item(82,1).PasteSpecial(-4122)
$ws.Cells.Item(82,1).Value = 80
$ws.Cells.Item(82,2).Value = 'Dollar'
$ws.Cells.Item(82,3).Value = '5,03'
$ws.Cells.Item(82,4).Value = '21:20'
$ws.Cells.Item(82,5).Value = ' quarta-feira, 5 de abril de 2023 '

$ws.Cells.Item(82,1).Copy()
$ws.Cells.Item(83,1).PasteSpecial(-4122)
$ws.Cells.Item(83,1).Value = 81
$ws.Cells.Item(83,2).Value = 'Euro'
$ws.Cells.Item(83,3).Value = '5,49'
$ws.Cells.Item(83,4).Value = '21:20'
$ws.Cells.Item(83,5).Value = ' quarta-feira, 5 de abril de 2023 '

$ws.Cells.Item(83,1).Copy()
$ws.Cells.Item(84,1).PasteSpecial(-4122)
$ws.Cells.Item(84,1).Value = 82
$ws.Cells.Item(84,2).Value = 'Iene'
$ws.Cells.Item(84,3).Value = '0,038'
$ws.Cells.Item(84,4).Value = '21:20'
$ws.Cells.Item(84,5).Value = ' quarta-feira, 5 de abril de 2023 '

$ws.Cells.Item(84,1).Copy()
$ws.Cells.Item(85,1).PasteSpecial(-4122)
$ws.Cells.Item(85,1).Value = 83
$ws.Cells.Item(85,2).Value = 'Yuan Chinês'
$ws.Cells.Item(85,3).Value = '0,73'
$ws.Cells.Item(85,4).Value = '21:20'
$ws.Cells.Item(85,5).Value = ' quarta-feira, 5 de abril de 2023 '

# Reset column C cells back to the default "Normal" style now that the text
# values are locked in, so they end up without any explicit number format.
foreach ($r in $cCells) { $ws.Cells.Item($r, 3).Style = "Normal" }

$excel.Application.CutCopyMode = $false
